# china ReEAP bar graph datawrapper
#
# The underlying data table that feeds the second ("RoEAP") bar chart is
# shifted one column to the left (L:N -> K:M) and the now-redundant
# helper column A (A18:A24) is cleared out. The chart's series formulas,
# the sheet selection, and the two chart frames' on-sheet positions are
# updated to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Shift the L:N (class / RoEAP_2017 / China_2017) helper table one
#    column to the left, into K:M, then clear the vacated A and N
#    columns. Cell-by-cell copy (rather than Range.Copy / ClearContents,
#    which mis-resolve overlapping ranges here) keeps this reliable.
# ---------------------------------------------------------------------
for ($col = 12; $col -le 14; $col++) {
    $destCol = $col - 1
    for ($row = 18; $row -le 24; $row++) {
        $srcCell = $ws.Cells.Item($row, $col)
        $val = $srcCell.Value2()
        $destCell = $ws.Cells.Item($row, $destCol)
        $destCell.Value = $val
    }
}

# Clear the now-vacated old column N (14) and helper column A.
for ($row = 18; $row -le 24; $row++) {
    $ws.Cells.Item($row, 14).ClearContents()
}
$ws.Range("A18:A24").ClearContents()

# ---------------------------------------------------------------------
# 2. Point the second bar chart's two series at the new K/L/M columns.
# ---------------------------------------------------------------------
$co2 = $ws.ChartObjects().Item(2)
$chart2 = $co2.Chart

$ser1 = $chart2.SeriesCollection().Item(1)
$ser1.Formula = "=SERIES('dreaded-bar-anime-data-for-use'!`$L`$18,'dreaded-bar-anime-data-for-use'!`$B`$19:`$B`$24,'dreaded-bar-anime-data-for-use'!`$L`$19:`$L`$24,1)"

$ser2 = $chart2.SeriesCollection().Item(2)
$ser2.Formula = "=SERIES('dreaded-bar-anime-data-for-use'!`$M`$18,'dreaded-bar-anime-data-for-use'!`$B`$19:`$B`$24,'dreaded-bar-anime-data-for-use'!`$M`$19:`$M`$24,2)"

# ---------------------------------------------------------------------
# 3. Reposition/resize both chart frames on the sheet.
# ---------------------------------------------------------------------
$co1 = $ws.ChartObjects().Item(1)
$co1.Left = 317.3681640625
$co1.Top = 236.5
$co1.Width = 343.6875
$co1.Height = 300.5

$co2.Left = 845.8056640625
$co2.Top = 238
$co2.Width = 343.6875
$co2.Height = 300.5

# ---------------------------------------------------------------------
# 4. Update the saved sheet selection.
# ---------------------------------------------------------------------
$ws.Range("K18:M24").Select() | Out-Null
